$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update workbook window view metadata ---
$win = $excel.ActiveWindow
$win.Left = 0
$win.Top = 500
$win.Width = 20480
$win.Height = 12680

# --- Update cell values for the Sample4/5/6/10 lane block rotation ---
$ws.Range("A6").Value = "Sample 10"
$ws.Range("B6").Value = "ERS4649260"
$ws.Range("C6").Value = "SAMEA6939367"
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = "male"
$ws.Range("U6").Value = "Sample 10"
$ws.Range("AS6").Value = "group 5"
$ws.Range("AU6").Value = "ERX4190630"
$ws.Range("AV6").Value = "S20271_156_S9_L001.bam"
$ws.Range("AW6").Value = "S20271_156_S9_L001.bam"
$ws.Range("AX6").Value = "S20271_156_S9_L001_R1_001.fastq.gz"
$ws.Range("AY6").Value = "S20271_156_S9_L001_R2_001.fastq.gz"
$ws.Range("BE6").Value = 100
$ws.Range("BI6").Value = "S20271_156_S9_L001"

$ws.Range("A7").Value = "Sample 10"
$ws.Range("B7").Value = "ERS4649260"
$ws.Range("C7").Value = "SAMEA6939367"
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = "male"
$ws.Range("U7").Value = "Sample 10"
$ws.Range("AS7").Value = "group 5"
$ws.Range("AU7").Value = "ERX4190630"
$ws.Range("AV7").Value = "S20271_156_S9_L002.bam"
$ws.Range("AW7").Value = "S20271_156_S9_L002.bam"
$ws.Range("AX7").Value = "S20271_156_S9_L002_R1_001.fastq.gz"
$ws.Range("AY7").Value = "S20271_156_S9_L002_R2_001.fastq.gz"
$ws.Range("BE7").Value = 100
$ws.Range("BI7").Value = "S20271_156_S9_L002"

$ws.Range("AV18").Value = "3391B_D50_1_S1_L001.bam"
$ws.Range("AW18").Value = "3391B_D50_1_S1_L001.bam"
$ws.Range("AX18").Value = "3391B_D50_1_S1_L001_R1_001.fastq.gz"
$ws.Range("AY18").Value = "3391B_D50_1_S1_L001_R2_001.fastq.gz"
$ws.Range("BI18").Value = "3391B_D50_1_S1_L001"

$ws.Range("AV19").Value = "3391B_D50_2_S1_L001.bam"
$ws.Range("AW19").Value = "3391B_D50_2_S1_L001.bam"
$ws.Range("AX19").Value = "3391B_D50_2_S1_L001_R1_001.fastq.gz"
$ws.Range("AY19").Value = "3391B_D50_2_S1_L001_R2_001.fastq.gz"
$ws.Range("BI19").Value = "3391B_D50_2_S1_L001"

$ws.Range("A20").Value = "Sample 4"
$ws.Range("B20").Value = "ERS4649261"
$ws.Range("C20").Value = "SAMEA6939368"
$ws.Range("E20").Value = 4
$ws.Range("G20").Value = "female"
$ws.Range("U20").Value = "Sample 4"
$ws.Range("AS20").Value = "group 6"
$ws.Range("AU20").Value = "ERX4190631"
$ws.Range("AV20").Value = "3391B_D50_3_S1_L001.bam"
$ws.Range("AW20").Value = "3391B_D50_3_S1_L001.bam"
$ws.Range("AX20").Value = "3391B_D50_3_S1_L001_R1_001.fastq.gz"
$ws.Range("AY20").Value = "3391B_D50_3_S1_L001_R2_001.fastq.gz"
$ws.Range("BI20").Value = "3391B_D50_3_S1_L001"

$ws.Range("A21").Value = "Sample 4"
$ws.Range("B21").Value = "ERS4649261"
$ws.Range("C21").Value = "SAMEA6939368"
$ws.Range("E21").Value = 4
$ws.Range("G21").Value = "female"
$ws.Range("U21").Value = "Sample 4"
$ws.Range("AS21").Value = "group 6"
$ws.Range("AU21").Value = "ERX4190631"
$ws.Range("AV21").Value = "3391B_D50_4_S1_L001.bam"
$ws.Range("AW21").Value = "3391B_D50_4_S1_L001.bam"
$ws.Range("AX21").Value = "3391B_D50_4_S1_L001_R1_001.fastq.gz"
$ws.Range("AY21").Value = "3391B_D50_4_S1_L001_R2_001.fastq.gz"
$ws.Range("BI21").Value = "3391B_D50_4_S1_L001"

$ws.Range("AV22").Value = "kolf2c1day50_1_S1_L001.bam"
$ws.Range("AW22").Value = "kolf2c1day50_1_S1_L001.bam"
$ws.Range("AX22").Value = "kolf2c1day50_1_S1_L001_R1_001.fastq.gz"
$ws.Range("AY22").Value = "kolf2c1day50_1_S1_L001_R2_001.fastq.gz"
$ws.Range("BI22").Value = "kolf2c1day50_1_S1_L001"

$ws.Range("AV23").Value = "kolf2c1day50_2_S1_L001.bam"
$ws.Range("AW23").Value = "kolf2c1day50_2_S1_L001.bam"
$ws.Range("AX23").Value = "kolf2c1day50_2_S1_L001_R1_001.fastq.gz"
$ws.Range("AY23").Value = "kolf2c1day50_2_S1_L001_R2_001.fastq.gz"
$ws.Range("BI23").Value = "kolf2c1day50_2_S1_L001"

$ws.Range("A24").Value = "Sample 5"
$ws.Range("B24").Value = "ERS4649262"
$ws.Range("C24").Value = "SAMEA6939369"
$ws.Range("E24").Value = 2
$ws.Range("U24").Value = "Sample 5"
$ws.Range("AS24").Value = "group 7"
$ws.Range("AU24").Value = "ERX4190632"
$ws.Range("AV24").Value = "kolf2c1day50_3_S1_L001.bam"
$ws.Range("AW24").Value = "kolf2c1day50_3_S1_L001.bam"
$ws.Range("AX24").Value = "kolf2c1day50_3_S1_L001_R1_001.fastq.gz"
$ws.Range("AY24").Value = "kolf2c1day50_3_S1_L001_R2_001.fastq.gz"
$ws.Range("BI24").Value = "kolf2c1day50_3_S1_L001"

$ws.Range("A25").Value = "Sample 5"
$ws.Range("B25").Value = "ERS4649262"
$ws.Range("C25").Value = "SAMEA6939369"
$ws.Range("E25").Value = 2
$ws.Range("U25").Value = "Sample 5"
$ws.Range("AS25").Value = "group 7"
$ws.Range("AU25").Value = "ERX4190632"
$ws.Range("AV25").Value = "kolf2c1day50_4_S1_L001.bam"
$ws.Range("AW25").Value = "kolf2c1day50_4_S1_L001.bam"
$ws.Range("AX25").Value = "kolf2c1day50_4_S1_L001_R1_001.fastq.gz"
$ws.Range("AY25").Value = "kolf2c1day50_4_S1_L001_R2_001.fastq.gz"
$ws.Range("BI25").Value = "kolf2c1day50_4_S1_L001"

$ws.Range("AV26").Value = "G1_MIF1D50_1_S1_L001.bam"
$ws.Range("AW26").Value = "G1_MIF1D50_1_S1_L001.bam"
$ws.Range("AX26").Value = "G1_MIF1D50_1_S1_L001_R1_001.fastq.gz"
$ws.Range("AY26").Value = "G1_MIF1D50_1_S1_L001_R2_001.fastq.gz"
$ws.Range("BI26").Value = "G1_MIF1D50_1_S1_L001"

$ws.Range("AV27").Value = "G1_MIF1D50_2_S1_L001.bam"
$ws.Range("AW27").Value = "G1_MIF1D50_2_S1_L001.bam"
$ws.Range("AX27").Value = "G1_MIF1D50_2_S1_L001_R1_001.fastq.gz"
$ws.Range("AY27").Value = "G1_MIF1D50_2_S1_L001_R2_001.fastq.gz"
$ws.Range("BI27").Value = "G1_MIF1D50_2_S1_L001"

$ws.Range("A28").Value = "Sample 6"
$ws.Range("B28").Value = "ERS4649263"
$ws.Range("C28").Value = "SAMEA6939370"
$ws.Range("U28").Value = "Sample 6"
$ws.Range("AS28").Value = "group 8"
$ws.Range("AU28").Value = "ERX4190633"
$ws.Range("AV28").Value = "G1_MIF1D50_3_S1_L001.bam"
$ws.Range("AW28").Value = "G1_MIF1D50_3_S1_L001.bam"
$ws.Range("AX28").Value = "G1_MIF1D50_3_S1_L001_R1_001.fastq.gz"
$ws.Range("AY28").Value = "G1_MIF1D50_3_S1_L001_R2_001.fastq.gz"
$ws.Range("BE28").Value = 50
$ws.Range("BI28").Value = "G1_MIF1D50_3_S1_L001"

$ws.Range("A29").Value = "Sample 6"
$ws.Range("B29").Value = "ERS4649263"
$ws.Range("C29").Value = "SAMEA6939370"
$ws.Range("U29").Value = "Sample 6"
$ws.Range("AS29").Value = "group 8"
$ws.Range("AU29").Value = "ERX4190633"
$ws.Range("AV29").Value = "G1_MIF1D50_4_S1_L001.bam"
$ws.Range("AW29").Value = "G1_MIF1D50_4_S1_L001.bam"
$ws.Range("AX29").Value = "G1_MIF1D50_4_S1_L001_R1_001.fastq.gz"
$ws.Range("AY29").Value = "G1_MIF1D50_4_S1_L001_R2_001.fastq.gz"
$ws.Range("BE29").Value = 50
$ws.Range("BI29").Value = "G1_MIF1D50_4_S1_L001"

# --- Rebuild hyperlinks on column BJ to follow the row rotation ---
$ws.Hyperlinks.Delete() | Out-Null
[void]$ws.Hyperlinks.Add($ws.Range("BJ2"), "https://www.ebi.ac.uk/ena/browser/view/ERR4229841", "", "", "https://www.ebi.ac.uk/ena/browser/view/ERR4229841")
[void]$ws.Hyperlinks.Add($ws.Range("BJ3"), "https://www.ebi.ac.uk/ena/browser/view/ERR4229842", "", "", "https://www.ebi.ac.uk/ena/browser/view/ERR4229842")
[void]$ws.Hyperlinks.Add($ws.Range("BJ4"), "https://www.ebi.ac.uk/ena/browser/view/ERR4229843", "", "", "https://www.ebi.ac.uk/ena/browser/view/ERR4229843")
[void]$ws.Hyperlinks.Add($ws.Range("BJ5"), "https://www.ebi.ac.uk/ena/browser/view/ERR4229844", "", "", "https://www.ebi.ac.uk/ena/browser/view/ERR4229844")
[void]$ws.Hyperlinks.Add($ws.Range("BJ6"), "https://www.ebi.ac.uk/ena/browser/view/ERR4229851", "", "", "https://www.ebi.ac.uk/ena/browser/view/ERR4229851")
[void]$ws.Hyperlinks.Add($ws.Range("BJ7"), "https://www.ebi.ac.uk/ena/browser/view/ERR4229852", "", "", "https://www.ebi.ac.uk/ena/browser/view/ERR4229852")
[void]$ws.Hyperlinks.Add($ws.Range("BJ8"), "https://www.ebi.ac.uk/ena/browser/view/ERR4229845", "", "", "https://www.ebi.ac.uk/ena/browser/view/ERR4229845")
[void]$ws.Hyperlinks.Add($ws.Range("BJ9"), "https://www.ebi.ac.uk/ena/browser/view/ERR4229846", "", "", "https://www.ebi.ac.uk/ena/browser/view/ERR4229846")
[void]$ws.Hyperlinks.Add($ws.Range("BJ10"), "https://www.ebi.ac.uk/ena/browser/view/ERR4229847", "", "", "https://www.ebi.ac.uk/ena/browser/view/ERR4229847")
[void]$ws.Hyperlinks.Add($ws.Range("BJ11"), "https://www.ebi.ac.uk/ena/browser/view/ERR4229848", "", "", "https://www.ebi.ac.uk/ena/browser/view/ERR4229848")
[void]$ws.Hyperlinks.Add($ws.Range("BJ12"), "https://www.ebi.ac.uk/ena/browser/view/ERR4229849", "", "", "https://www.ebi.ac.uk/ena/browser/view/ERR4229849")
[void]$ws.Hyperlinks.Add($ws.Range("BJ13"), "https://www.ebi.ac.uk/ena/browser/view/ERR4229850", "", "", "https://www.ebi.ac.uk/ena/browser/view/ERR4229850")
[void]$ws.Hyperlinks.Add($ws.Range("BJ14"), "https://www.ebi.ac.uk/ena/browser/view/ERR4229837", "", "", "https://www.ebi.ac.uk/ena/browser/view/ERR4229837")
[void]$ws.Hyperlinks.Add($ws.Range("BJ15"), "https://www.ebi.ac.uk/ena/browser/view/ERR4229838", "", "", "https://www.ebi.ac.uk/ena/browser/view/ERR4229838")
[void]$ws.Hyperlinks.Add($ws.Range("BJ16"), "https://www.ebi.ac.uk/ena/browser/view/ERR4229839", "", "", "https://www.ebi.ac.uk/ena/browser/view/ERR4229839")
[void]$ws.Hyperlinks.Add($ws.Range("BJ17"), "https://www.ebi.ac.uk/ena/browser/view/ERR4229840", "", "", "https://www.ebi.ac.uk/ena/browser/view/ERR4229840")
[void]$ws.Hyperlinks.Add($ws.Range("BJ18"), "https://www.ebi.ac.uk/ena/browser/view/ERR4229853", "", "", "https://www.ebi.ac.uk/ena/browser/view/ERR4229853")
[void]$ws.Hyperlinks.Add($ws.Range("BJ19"), "https://www.ebi.ac.uk/ena/browser/view/ERR4229854", "", "", "https://www.ebi.ac.uk/ena/browser/view/ERR4229854")
[void]$ws.Hyperlinks.Add($ws.Range("BJ20"), "https://www.ebi.ac.uk/ena/browser/view/ERR4229855", "", "", "https://www.ebi.ac.uk/ena/browser/view/ERR4229855")
[void]$ws.Hyperlinks.Add($ws.Range("BJ21"), "https://www.ebi.ac.uk/ena/browser/view/ERR4229856", "", "", "https://www.ebi.ac.uk/ena/browser/view/ERR4229856")
[void]$ws.Hyperlinks.Add($ws.Range("BJ22"), "https://www.ebi.ac.uk/ena/browser/view/ERR4229857", "", "", "https://www.ebi.ac.uk/ena/browser/view/ERR4229857")
[void]$ws.Hyperlinks.Add($ws.Range("BJ23"), "https://www.ebi.ac.uk/ena/browser/view/ERR4229858", "", "", "https://www.ebi.ac.uk/ena/browser/view/ERR4229858")
[void]$ws.Hyperlinks.Add($ws.Range("BJ24"), "https://www.ebi.ac.uk/ena/browser/view/ERR4229859", "", "", "https://www.ebi.ac.uk/ena/browser/view/ERR4229859")
[void]$ws.Hyperlinks.Add($ws.Range("BJ25"), "https://www.ebi.ac.uk/ena/browser/view/ERR4229860", "", "", "https://www.ebi.ac.uk/ena/browser/view/ERR4229860")
[void]$ws.Hyperlinks.Add($ws.Range("BJ26"), "https://www.ebi.ac.uk/ena/browser/view/ERR4229861", "", "", "https://www.ebi.ac.uk/ena/browser/view/ERR4229861")
[void]$ws.Hyperlinks.Add($ws.Range("BJ27"), "https://www.ebi.ac.uk/ena/browser/view/ERR4229862", "", "", "https://www.ebi.ac.uk/ena/browser/view/ERR4229862")
[void]$ws.Hyperlinks.Add($ws.Range("BJ28"), "https://www.ebi.ac.uk/ena/browser/view/ERR4229863", "", "", "https://www.ebi.ac.uk/ena/browser/view/ERR4229863")
[void]$ws.Hyperlinks.Add($ws.Range("BJ29"), "https://www.ebi.ac.uk/ena/browser/view/ERR4229864", "", "", "https://www.ebi.ac.uk/ena/browser/view/ERR4229864")
[void]$ws.Hyperlinks.Add($ws.Range("BJ30"), "https://www.ebi.ac.uk/ena/browser/view/ERR4229865", "", "", "https://www.ebi.ac.uk/ena/browser/view/ERR4229865")
[void]$ws.Hyperlinks.Add($ws.Range("BJ31"), "https://www.ebi.ac.uk/ena/browser/view/ERR4229866", "", "", "https://www.ebi.ac.uk/ena/browser/view/ERR4229866")
[void]$ws.Hyperlinks.Add($ws.Range("BJ32"), "https://www.ebi.ac.uk/ena/browser/view/ERR4229867", "", "", "https://www.ebi.ac.uk/ena/browser/view/ERR4229867")
[void]$ws.Hyperlinks.Add($ws.Range("BJ33"), "https://www.ebi.ac.uk/ena/browser/view/ERR4229868", "", "", "https://www.ebi.ac.uk/ena/browser/view/ERR4229868")
[void]$ws.Hyperlinks.Add($ws.Range("BJ34"), "https://www.ebi.ac.uk/ena/browser/view/ERR4229869", "", "", "https://www.ebi.ac.uk/ena/browser/view/ERR4229869")
[void]$ws.Hyperlinks.Add($ws.Range("BJ35"), "https://www.ebi.ac.uk/ena/browser/view/ERR4229870", "", "", "https://www.ebi.ac.uk/ena/browser/view/ERR4229870")

# --- Update the active sheet view selection/navigation state ---
[void]$ws.Range("AX34").Select()
